$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(86, 8).Value = 6123
$ws_ALC.Cells.Item(86, 9).Value = 5984.5
$ws_ALC.Cells.Item(86, 10).Value = 6400
$ws_ALC.Cells.Item(86, 11).Value = 5984.5
$ws_ALC.Cells.Item(86, 12).Value = 6400
$ws_ALC.Cells.Item(86, 13).Value = -4861.5
$ws_ALC.Cells.Item(86, 14).Value = -8646
$ws_ALC.Cells.Item(89, 8).Value = 6123
$ws_ALC.Cells.Item(89, 9).Value = 5984.5
$ws_ALC.Cells.Item(89, 10).Value = 6400
$ws_ALC.Cells.Item(89, 11).Value = 29922.5
$ws_ALC.Cells.Item(89, 12).Value = 32000
$ws_ALC.Cells.Item(89, 13).Value = -24306.5
$ws_ALC.Cells.Item(89, 14).Value = -43232
$ws_ALC.Cells.Item(96, 8).Value = 556.8182
$ws_ALC.Cells.Item(96, 9).Value = 479
$ws_ALC.Cells.Item(96, 10).Value = 764.3333
$ws_ALC.Cells.Item(96, 11).Value = 1437
$ws_ALC.Cells.Item(96, 12).Value = 2292.9999
$ws_ALC.Cells.Item(96, 13).Value = -64
$ws_ALC.Cells.Item(96, 14).Value = -5038.9999
$ws_ALC.Cells.Item(113, 8).Value = 6251.78
$ws_ALC.Cells.Item(113, 9).Value = 6237.0884
$ws_ALC.Cells.Item(113, 11).Value = 6237.0884
$ws_ALC.Cells.Item(113, 13).Value = -2983.0884
$ws_ALC.Cells.Item(116, 8).Value = 8921.375
$ws_ALC.Cells.Item(116, 9).Value = 4899.5
$ws_ALC.Cells.Item(116, 11).Value = 4899.5
$ws_ALC.Cells.Item(116, 13).Value = -1457.5
$ws_ALC.Cells.Item(132, 8).Value = 2591.913
$ws_ALC.Cells.Item(132, 9).Value = 2338.7856
$ws_ALC.Cells.Item(132, 11).Value = 7016.3568
$ws_ALC.Cells.Item(132, 13).Value = -4486.3568
$ws_ALC.Cells.Item(137, 8).Value = 502053.2
$ws_ALC.Cells.Item(137, 9).Value = 770909
$ws_ALC.Cells.Item(137, 10).Value = 2749.5715
$ws_ALC.Cells.Item(137, 11).Value = 2312727
$ws_ALC.Cells.Item(137, 12).Value = 8248.7145
$ws_ALC.Cells.Item(137, 13).Value = -2310177
$ws_ALC.Cells.Item(137, 14).Value = -13348.7145
$ws_ALC.Cells.Item(138, 8).Value = 4150.0625
$ws_ALC.Cells.Item(138, 9).Value = 4586.2
$ws_ALC.Cells.Item(138, 10).Value = 3765.2354
$ws_ALC.Cells.Item(138, 11).Value = 13758.6
$ws_ALC.Cells.Item(138, 12).Value = 11295.7062
$ws_ALC.Cells.Item(138, 13).Value = -8618.599999999999
$ws_ALC.Cells.Item(138, 14).Value = -21575.7062
$ws_ALC.Cells.Item(141, 8).Value = 4670.357
$ws_ALC.Cells.Item(141, 9).Value = 2736.2917
$ws_ALC.Cells.Item(141, 10).Value = 16274.75
$ws_ALC.Cells.Item(141, 11).Value = 8208.875100000001
$ws_ALC.Cells.Item(141, 12).Value = 48824.25
$ws_ALC.Cells.Item(141, 13).Value = -3028.875100000001
$ws_ALC.Cells.Item(141, 14).Value = -59184.25
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(74, 8).Value = 3620.28
$ws_ARM.Cells.Item(74, 9).Value = 1229.2354
$ws_ARM.Cells.Item(74, 10).Value = 8701.25
$ws_ARM.Cells.Item(74, 11).Value = 1229.2354
$ws_ARM.Cells.Item(74, 12).Value = 8701.25
$ws_ARM.Cells.Item(74, 13).Value = -355.2354
$ws_ARM.Cells.Item(74, 14).Value = -10449.25
$ws_ARM.Cells.Item(77, 8).Value = 3620.28
$ws_ARM.Cells.Item(77, 9).Value = 1229.2354
$ws_ARM.Cells.Item(77, 10).Value = 8701.25
$ws_ARM.Cells.Item(77, 11).Value = 6146.177
$ws_ARM.Cells.Item(77, 12).Value = 43506.25
$ws_ARM.Cells.Item(77, 13).Value = -1778.177
$ws_ARM.Cells.Item(77, 14).Value = -52242.25
$ws_ARM.Cells.Item(110, 8).Value = 80559.28999999999
$ws_ARM.Cells.Item(110, 9).Value = 93584.25
$ws_ARM.Cells.Item(110, 11).Value = 93584.25
$ws_ARM.Cells.Item(110, 13).Value = -91539.25
$ws_ARM.Cells.Item(132, 8).Value = 1885.5
$ws_ARM.Cells.Item(132, 9).Value = 1426.2632
$ws_ARM.Cells.Item(132, 11).Value = 4278.7896
$ws_ARM.Cells.Item(132, 13).Value = -1748.7896
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(21, 8).Value = 17840.25
$ws_BSM.Cells.Item(21, 10).Value = 17840.25
$ws_BSM.Cells.Item(21, 12).Value = 17840.25
$ws_BSM.Cells.Item(21, 14).Value = -18312.25
$ws_BSM.Cells.Item(86, 8).Value = 11274.77
$ws_BSM.Cells.Item(86, 10).Value = 17321.875
$ws_BSM.Cells.Item(86, 12).Value = 17321.875
$ws_BSM.Cells.Item(86, 14).Value = -19567.875
$ws_BSM.Cells.Item(89, 8).Value = 11274.77
$ws_BSM.Cells.Item(89, 10).Value = 17321.875
$ws_BSM.Cells.Item(89, 12).Value = 86609.375
$ws_BSM.Cells.Item(89, 14).Value = -97841.375
$ws_BSM.Cells.Item(94, 8).Value = 2945586
$ws_BSM.Cells.Item(94, 9).Value = 4168099.8
$ws_BSM.Cells.Item(94, 11).Value = 4168099.8
$ws_BSM.Cells.Item(94, 13).Value = -4167648.8
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Cells.Item(22, 8).Value = 720.1667
$ws_CRP.Cells.Item(22, 10).Value = 832.75
$ws_CRP.Cells.Item(22, 12).Value = 832.75
$ws_CRP.Cells.Item(22, 14).Value = -1532.75
$ws_CRP.Cells.Item(74, 8).Value = 28892
$ws_CRP.Cells.Item(74, 10).Value = 35094.332
$ws_CRP.Cells.Item(74, 12).Value = 35094.332
$ws_CRP.Cells.Item(74, 14).Value = -36842.332
$ws_CRP.Cells.Item(77, 8).Value = 28892
$ws_CRP.Cells.Item(77, 10).Value = 35094.332
$ws_CRP.Cells.Item(77, 12).Value = 105282.996
$ws_CRP.Cells.Item(77, 14).Value = -114018.996
$ws_CRP.Cells.Item(99, 8).Value = 2731.6
$ws_CRP.Cells.Item(99, 9).Value = 2512.9333
$ws_CRP.Cells.Item(99, 11).Value = 2512.9333
$ws_CRP.Cells.Item(99, 13).Value = -1014.9333
$ws_CRP.Cells.Item(126, 8).Value = 2731.6
$ws_CRP.Cells.Item(126, 9).Value = 2512.9333
$ws_CRP.Cells.Item(126, 11).Value = 7538.7999
$ws_CRP.Cells.Item(126, 13).Value = -5068.7999
$ws_CRP.Cells.Item(132, 8).Value = 1991.0588
$ws_CRP.Cells.Item(132, 9).Value = 1896.7812
$ws_CRP.Cells.Item(132, 11).Value = 5690.3436
$ws_CRP.Cells.Item(132, 13).Value = -3160.3436
$ws_CRP.Cells.Item(134, 8).Value = 3503.9092
$ws_CRP.Cells.Item(134, 9).Value = 2549.6667
$ws_CRP.Cells.Item(134, 11).Value = 7649.000100000001
$ws_CRP.Cells.Item(134, 13).Value = -5114.000100000001
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(113, 8).Value = 3414.111
$ws_GSM.Cells.Item(113, 9).Value = 2174.75
$ws_GSM.Cells.Item(113, 11).Value = 2174.75
$ws_GSM.Cells.Item(113, 13).Value = -4.75
$ws_GSM.Cells.Item(122, 8).Value = 3669.889
$ws_GSM.Cells.Item(122, 9).Value = 3607.182
$ws_GSM.Cells.Item(122, 10).Value = 3768.4285
$ws_GSM.Cells.Item(122, 11).Value = 10821.546
$ws_GSM.Cells.Item(122, 12).Value = 11305.2855
$ws_GSM.Cells.Item(122, 13).Value = -8371.545999999998
$ws_GSM.Cells.Item(122, 14).Value = -16205.2855
$ws_GSM.Cells.Item(132, 8).Value = 4290.1934
$ws_GSM.Cells.Item(132, 9).Value = 4534.115
$ws_GSM.Cells.Item(132, 10).Value = 3021.8
$ws_GSM.Cells.Item(132, 11).Value = 13602.345
$ws_GSM.Cells.Item(132, 12).Value = 9065.400000000001
$ws_GSM.Cells.Item(132, 13).Value = -11072.345
$ws_GSM.Cells.Item(132, 14).Value = -14125.4
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(40, 8).Value = 1652.15
$ws_LTW.Cells.Item(40, 9).Value = 1480.2778
$ws_LTW.Cells.Item(40, 11).Value = 1480.2778
$ws_LTW.Cells.Item(40, 13).Value = -1344.2778
$ws_LTW.Cells.Item(55, 8).Value = 7143269.5
$ws_LTW.Cells.Item(55, 9).Value = 11111345
$ws_LTW.Cells.Item(55, 11).Value = 11111345
$ws_LTW.Cells.Item(55, 13).Value = -11111172
$ws_LTW.Cells.Item(95, 8).Value = 20666.666
$ws_LTW.Cells.Item(95, 10).Value = 20666.666
$ws_LTW.Cells.Item(95, 12).Value = 20666.666
$ws_LTW.Cells.Item(95, 14).Value = -26158.666
$ws_LTW.Cells.Item(100, 8).Value = 859.875
$ws_LTW.Cells.Item(100, 9).Value = 697
$ws_LTW.Cells.Item(100, 11).Value = 697
$ws_LTW.Cells.Item(100, 13).Value = -156
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Cells.Item(107, 8).Value = 7677.933
$ws_WVR.Cells.Item(107, 9).Value = 13397
$ws_WVR.Cells.Item(107, 11).Value = 40191
$ws_WVR.Cells.Item(107, 13).Value = -38271
$ws_WVR.Cells.Item(132, 8).Value = 1478.5652
$ws_WVR.Cells.Item(132, 9).Value = 1383.6904
$ws_WVR.Cells.Item(132, 11).Value = 4151.0712
$ws_WVR.Cells.Item(132, 13).Value = -1621.0712
$ws_WVR.Cells.Item(136, 8).Value = 9290.612999999999
$ws_WVR.Cells.Item(136, 9).Value = 8433.634
$ws_WVR.Cells.Item(136, 11).Value = 25300.902
$ws_WVR.Cells.Item(136, 13).Value = -22750.902
